$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "Johns Hopkins 데이터 과학 석사 프로그램 소개"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/johns-hopkins-msds/#utm_source=rss&utm_medium=rss&utm_campaign=johns-hopkins-msds"

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D27").Value = "ACL 2022 Review"
$ws.Range("E27").Value = "https://blog.pingpong.us/acl2022-review/"

$ws.Range("D28").Value = "[테스트 영상] Simple Online and Realtime Tracking"
$ws.Range("E28").Value = "https://ropiens.tistory.com/189"

$ws.Range("D51").Value = "[윈도우11] 디스플레이 해상도 변경이 안 될 때 조치 방법"
$ws.Range("E51").Value = "https://bskyvision.com/1295"

$ws.Range("D52").Value = "숨은 DS"
